$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.054431
$ws.Range("I2").Value = 0.06647005666275889
$ws.Range("J2").Value = 0.06647005666275889
$ws.Range("M2").Value = 1.037532
$ws.Range("N2").Value = 3.112596
$ws.Range("O2").Value = 0.04166450179684251
$ws.Range("P2").Value = 0.0439159257402554
$ws.Range("Q2").Value = 0.018824634764
$ws.Range("R2").Value = 0.169421712876
$ws.Range("S2").Value = 0.002769441795261741
$ws.Range("T2").Value = 0.002919094072352288
$ws.Range("H3").Value = 0.054431
$ws.Range("I3").Value = 0.06647005666275889
$ws.Range("J3").Value = 0.06647005666275889
$ws.Range("O3").Value = 0.2561129158441639
$ws.Range("P3").Value = 0.2699524849277078
$ws.Range("S3").Value = 0.01702384002822597
$ws.Range("T3").Value = 0.0179437569693973
$ws.Range("H4").Value = 0.054431
$ws.Range("I4").Value = 0.06647005666275889
$ws.Range("J4").Value = 0.06647005666275889
$ws.Range("M4").Value = 6.239319333333333
$ws.Range("N4").Value = 18.717958
$ws.Range("O4").Value = 0.2505543265891952
$ws.Range("P4").Value = 0.2640935262839185
$ws.Range("Q4").Value = 0.1132041302108889
$ws.Range("R4").Value = 1.018837171898
$ws.Range("S4").Value = 0.0166543602854832
$ws.Range("T4").Value = 0.01755431165635987
$ws.Range("H5").Value = 0.054431
$ws.Range("I5").Value = 0.06647005666275889
$ws.Range("J5").Value = 0.06647005666275889
$ws.Range("M5").Value = 3.8299385
$ws.Range("N5").Value = 7.659877
$ws.Range("O5").Value = 0.1538000558200097
$ws.Range("P5").Value = 0.1080739644693659
$ws.Range("Q5").Value = 0.06948912749783333
$ws.Range("R5").Value = 0.416934764987
$ws.Range("S5").Value = 0.01022309842509152
$ws.Range("T5").Value = 0.007183682542047739
$ws.Range("H6").Value = 0.054431
$ws.Range("I6").Value = 0.06647005666275889
$ws.Range("J6").Value = 0.06647005666275889
$ws.Range("M6").Value = 7.417532333333334
$ws.Range("N6").Value = 22.252597
$ws.Range("O6").Value = 0.2978681999497886
$ws.Range("P6").Value = 0.3139640985787523
$ws.Range("Q6").Value = 0.1345812341452222
$ws.Range("R6").Value = 1.211231107307
$ws.Range("S6").Value = 0.01979931612869645
$ws.Range("T6").Value = 0.02086921142260168
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.10505
$ws.Range("H7").Value = 0.31515
$ws.Range("I7").Value = 0.3848549237983587
$ws.Range("J7").Value = 0.3848549237983587
$ws.Range("M7").Value = 1.037532
$ws.Range("N7").Value = 3.112596
$ws.Range("O7").Value = 0.04166450179684251
$ws.Range("P7").Value = 0.0439159257402554
$ws.Range("Q7").Value = 0.1089927366
$ws.Range("R7").Value = 0.9809346294
$ws.Range("S7").Value = 0.0160347886641204
$ws.Range("T7").Value = 0.01690126025430037
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.10505
$ws.Range("H8").Value = 0.31515
$ws.Range("I8").Value = 0.3848549237983587
$ws.Range("J8").Value = 0.3848549237983587
$ws.Range("O8").Value = 0.2561129158441639
$ws.Range("P8").Value = 0.2699524849277078
$ws.Range("Q8").Value = 0.6699815519833332
$ws.Range("R8").Value = 6.029833967849998
$ws.Range("S8").Value = 0.09856631671098116
$ws.Range("T8").Value = 0.1038925430160306
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.10505
$ws.Range("H9").Value = 0.31515
$ws.Range("I9").Value = 0.3848549237983587
$ws.Range("J9").Value = 0.3848549237983587
$ws.Range("M9").Value = 6.239319333333333
$ws.Range("N9").Value = 18.717958
$ws.Range("O9").Value = 0.2505543265891952
$ws.Range("P9").Value = 0.2640935262839185
$ws.Range("Q9").Value = 0.6554404959666666
$ws.Range("R9").Value = 5.8989644637
$ws.Range("S9").Value = 0.0964270662668338
$ws.Range("T9").Value = 0.1016376939336373
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.10505
$ws.Range("H10").Value = 0.31515
$ws.Range("I10").Value = 0.3848549237983587
$ws.Range("J10").Value = 0.3848549237983587
$ws.Range("M10").Value = 3.8299385
$ws.Range("N10").Value = 7.659877
$ws.Range("O10").Value = 0.1538000558200097
$ws.Range("P10").Value = 0.1080739644693659
$ws.Range("Q10").Value = 0.4023350394249999
$ws.Range("R10").Value = 2.41401023655
$ws.Range("S10").Value = 0.05919070876279315
$ws.Range("T10").Value = 0.04159279736044432
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.10505
$ws.Range("H11").Value = 0.31515
$ws.Range("I11").Value = 0.3848549237983587
$ws.Range("J11").Value = 0.3848549237983587
$ws.Range("M11").Value = 7.417532333333334
$ws.Range("N11").Value = 22.252597
$ws.Range("O11").Value = 0.2978681999497886
$ws.Range("P11").Value = 0.3139640985787523
$ws.Range("Q11").Value = 0.7792117716166667
$ws.Range("R11").Value = 7.01290594455
$ws.Range("S11").Value = 0.1146360433936302
$ws.Range("T11").Value = 0.1208306292339461
$ws.Range("G12").Value = 0.1497663333333333
$ws.Range("H12").Value = 0.449299
$ws.Range("I12").Value = 0.5486750195388824
$ws.Range("J12").Value = 0.5486750195388823
$ws.Range("M12").Value = 1.037532
$ws.Range("N12").Value = 3.112596
$ws.Range("O12").Value = 0.04166450179684251
$ws.Range("P12").Value = 0.0439159257402554
$ws.Range("Q12").Value = 0.155387363356
$ws.Range("R12").Value = 1.398486270204
$ws.Range("S12").Value = 0.02286027133746037
$ws.Range("T12").Value = 0.02409557141360274
$ws.Range("G13").Value = 0.1497663333333333
$ws.Range("H13").Value = 0.449299
$ws.Range("I13").Value = 0.5486750195388824
$ws.Range("J13").Value = 0.5486750195388823
$ws.Range("O13").Value = 0.2561129158441639
$ws.Range("P13").Value = 0.2699524849277078
$ws.Range("Q13").Value = 0.9551706848312221
$ws.Range("R13").Value = 8.596536163480998
$ws.Range("S13").Value = 0.1405227591049568
$ws.Range("T13").Value = 0.1481161849422799
$ws.Range("G14").Value = 0.1497663333333333
$ws.Range("H14").Value = 0.449299
$ws.Range("I14").Value = 0.5486750195388824
$ws.Range("J14").Value = 0.5486750195388823
$ws.Range("M14").Value = 6.239319333333333
$ws.Range("N14").Value = 18.717958
$ws.Range("O14").Value = 0.2505543265891952
$ws.Range("P14").Value = 0.2640935262839185
$ws.Range("Q14").Value = 0.9344399790491111
$ws.Range("R14").Value = 8.409959811442
$ws.Range("S14").Value = 0.1374729000368782
$ws.Range("T14").Value = 0.1449015206939213
$ws.Range("G15").Value = 0.1497663333333333
$ws.Range("H15").Value = 0.449299
$ws.Range("I15").Value = 0.5486750195388824
$ws.Range("J15").Value = 0.5486750195388823
$ws.Range("M15").Value = 3.8299385
$ws.Range("N15").Value = 7.659877
$ws.Range("O15").Value = 0.1538000558200097
$ws.Range("P15").Value = 0.1080739644693659
$ws.Range("Q15").Value = 0.5735958460371666
$ws.Range("R15").Value = 3.441575076223
$ws.Range("S15").Value = 0.08438624863212503
$ws.Range("T15").Value = 0.05929748456687379
$ws.Range("G16").Value = 0.1497663333333333
$ws.Range("H16").Value = 0.449299
$ws.Range("I16").Value = 0.5486750195388824
$ws.Range("J16").Value = 0.5486750195388823
$ws.Range("M16").Value = 7.417532333333334
$ws.Range("N16").Value = 22.252597
$ws.Range("O16").Value = 0.2978681999497886
$ws.Range("P16").Value = 0.3139640985787523
$ws.Range("Q16").Value = 1.110896619944778
$ws.Range("R16").Value = 9.998069579503001
$ws.Range("S16").Value = 0.163432840427462
$ws.Range("T16").Value = 0.1722642579222045
